$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 293, shifting existing rows 293-410 down to 294-411.
$ws.Rows.Item(293).Insert()

# Populate the newly inserted row 293 with the new data record.
$ws.Range("A293").Value = 10
$ws.Range("B293").Value = "Vega Modelo de Temuco"
$ws.Range("C293").Value = "La Araucanía"
$ws.Range("D293").Value = 45141
$ws.Range("E293").Value = 9
$ws.Range("F293").Value = 100112039
$ws.Range("G293").Value = "Ciboulette"
$ws.Range("H293").Value = "Sin especificar"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 40
$ws.Range("K293").Value = 7000
$ws.Range("L293").Value = 7000
$ws.Range("M293").Value = 7000
$ws.Range("N293").Value = "$/docena de atados"
$ws.Range("O293").Value = "Provincia de Cautín"
$ws.Range("P293").Value = 2333
$ws.Range("Q293").Value = 3
$ws.Range("R293").Value = "Hortaliza"
